# Applies the "第十周周一" (week 10, Monday) progress block to the
# project plan sheet, plus small fixes to the preceding (week 9 Thu) block.
#
# Style legend used throughout this sheet (cf. existing rows 1/8/15/...):
#   style "title"  : thin border all sides, bold 10pt font, center+center align
#                    -> used for the merged "日期：..." banner row
#   style "header" : thin border all sides, bold 10pt font, default (vertical-center) align
#                    -> used for the 组员/计划内容/完成情况/备注 header row
#   style "plain"  : thin border all sides, regular font, default align
#                    -> used for plain data cells
#   style "left"   : thin border all sides, regular font, left+center align
#                    -> used for the "总结：" summary cell
#   style "pct"    : thin border all sides, regular font, percentage number format
#                    -> used for one "完成情况" cell per block (kept even when it
#                       ends up holding text instead of a number)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TitleStyle($rng) {
    $rng.Borders.LineStyle = 1
    $rng.Font.Bold = $true
    $rng.Font.Size = 10
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

function Set-HeaderStyle($rng) {
    $rng.Borders.LineStyle = 1
    $rng.Font.Bold = $true
    $rng.Font.Size = 10
}

function Set-PlainStyle($rng) {
    $rng.Borders.LineStyle = 1
}

function Set-LeftStyle($rng) {
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4131
}

function Set-PercentStyle($rng) {
    $rng.Borders.LineStyle = 1
    $rng.NumberFormat = "0%"
}

# ---------------------------------------------------------------------
# 1) Fixes inside the existing "第九周周四" block (rows 45-49)
# ---------------------------------------------------------------------
$ws.Range("B46").Value = "阅读习题题目与答案"
$ws.Range("C46").Value = "完成"

# ---------------------------------------------------------------------
# 2) Correct the block-43 banner date (周三 -> 周四)
# ---------------------------------------------------------------------
$ws.Range("A50").Value = "日期：2018.10.31 第九周周四"

# ---------------------------------------------------------------------
# 3) Append the new "第十周周一" block in rows 57-63
# ---------------------------------------------------------------------
$ws.Range("A57").Value = "日期：2018.11.5 第十周周一"
$ws.Range("A57:D57").Merge()
Set-TitleStyle($ws.Range("A57:D57"))

$ws.Range("A58").Value = "组员"
$ws.Range("B58").Value = "计划内容"
$ws.Range("C58").Value = "完成情况"
$ws.Range("D58").Value = "备注"
Set-HeaderStyle($ws.Range("A58:D58"))

$ws.Range("A59").Value = "苏立明"
$ws.Range("B59").Value = "添加搜索功能"
$ws.Range("C59").Value = "进行中"
Set-PlainStyle($ws.Range("A59:D59"))

$ws.Range("A60").Value = "何舒静"
$ws.Range("B60").Value = "添加登录验证"
$ws.Range("C60").Value = "进行中"
Set-PlainStyle($ws.Range("A60:D60"))
Set-PercentStyle($ws.Range("C60"))

$ws.Range("A61").Value = "郑瑞贤"
$ws.Range("B61").Value = "修改考试界面"
$ws.Range("C61").Value = "进行中"
Set-PlainStyle($ws.Range("A61:D61"))

$ws.Range("A62").Value = "总结："
Set-LeftStyle($ws.Range("A62:D62"))
Set-PlainStyle($ws.Range("A63:D63"))
$ws.Range("A62:D63").Merge()

# ---------------------------------------------------------------------
# 4) Update the saved view state (matches where the author was scrolled to)
# ---------------------------------------------------------------------
$ws.Range("A50:D50").Select()
$excel.ActiveWindow.ScrollRow = 25
